$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-296: update the Date_seconds (column A) and Date (column K) values
for ($r = 2; $r -le 296; $r++) {
    $ws.Cells.Item($r, 1).Value = 1529520284
    $ws.Cells.Item($r, 11).Value = "Wed_Jun_20_14:44:44_EDT_2018"
}

# Row 297: distinct updates
$ws.Cells.Item(297, 1).Value = 1529521748
$ws.Cells.Item(297, 2).Value = 1
$ws.Cells.Item(297, 4).Value = " Manufacturing "
$ws.Cells.Item(297, 11).Value = " Wed_Jun_20_15:09:08_EDT_2018"
